$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.03827233618454998
$ws.Range("H2").Value = 0.03247769798436308
$ws.Range("I2").Value = 0.8877442450792397
$ws.Range("G3").Value = 0.01722336798284841
$ws.Range("H3").Value = 0.01036823397201161
$ws.Range("I3").Value = 0.97777193124281
$ws.Range("G4").Value = 0.1100806331996572
$ws.Range("H4").Value = 0.0894070994215852
$ws.Range("I4").Value = 0.6729155746410738
$ws.Range("G5").Value = 0.02896702165290135
$ws.Range("H5").Value = 0.023198438221637
$ws.Range("I5").Value = 0.8948687861483507
$ws.Range("G8").Value = 0.04881769538406525
$ws.Range("H8").Value = 0.03988679807555299
$ws.Range("I8").Value = 0.8761875282402057
$ws.Range("G9").Value = 0.0153954658819749
$ws.Range("H9").Value = 0.01277548275624916
$ws.Range("I9").Value = 0.9778221689690755
$ws.Range("G10").Value = 0.01953871577594613
$ws.Range("H10").Value = 0.01525361413894663
$ws.Range("I10").Value = 0.963392733986809
$ws.Range("G11").Value = 0.01527694327980593
$ws.Range("H11").Value = 0.01116844329407299
$ws.Range("I11").Value = 0.9819019886256762
$ws.Range("G12").Value = 0.009398219727310895
$ws.Range("H12").Value = 0.007985011453821937
$ws.Range("I12").Value = 0.9908878264805299
$ws.Range("G13").Value = 0.043011340495216
$ws.Range("H13").Value = 0.03285187026852762
$ws.Range("I13").Value = 0.87144661793476
$ws.Range("G15").Value = 0.0114329995726591
$ws.Range("H15").Value = 0.009031354478747148
$ws.Range("I15").Value = 0.988939352242425
$ws.Range("G16").Value = 0.03402469959302017
$ws.Range("H16").Value = 0.02759517351137924
$ws.Range("I16").Value = 0.8919195594780989
$ws.Range("G18").Value = 0.009365056076055239
$ws.Range("H18").Value = 0.007805381617066186
$ws.Range("I18").Value = 0.9916206177117171
$ws.Range("G19").Value = 0.06893383870843871
$ws.Range("H19").Value = 0.05893689614694798
$ws.Range("I19").Value = 0.6596614666864518
$ws.Range("G20").Value = 0.08758357601932995
$ws.Range("H20").Value = 0.08193394843334012
$ws.Range("I20").Value = 0.6533546374326038
$ws.Range("G21").Value = 0.009958599028124186
$ws.Range("H21").Value = 0.008979613300267827
$ws.Range("I21").Value = 0.9932731385953062
$ws.Range("G25").Value = 0.06983743145843234
$ws.Range("H25").Value = 0.05720128196113361
$ws.Range("I25").Value = 0.8170976857267184
$ws.Range("G28").Value = 0.004053788344421826
$ws.Range("H28").Value = 0.003372044329074342
$ws.Range("I28").Value = 0.9991331937491382
$ws.Range("G30").Value = 0.001259307887700345
$ws.Range("H30").Value = 0.001259307887700345
$ws.Range("I30").Value = 0
$ws.Range("G31").Value = 0.03102470131072205
$ws.Range("H31").Value = 0.02507282769126183
$ws.Range("I31").Value = 0.9471220255892986
$ws.Range("G32").Value = 0.1105565951751388
$ws.Range("H32").Value = 0.07936349729831112
$ws.Range("I32").Value = 0.5440696468150541
$ws.Range("G33").Value = 0.04551055279069954
$ws.Range("H33").Value = 0.03772020820834515
$ws.Range("I33").Value = 0.8596485676655139
$ws.Range("G34").Value = 0.03945000544631413
$ws.Range("H34").Value = 0.03383294239672407
$ws.Range("I34").Value = 0.8714344694802804
$ws.Range("G35").Value = 0.04587632670855927
$ws.Range("H35").Value = 0.0330300477873342
$ws.Range("I35").Value = 0.8329578218690425
$ws.Range("G36").Value = 0.01885899505183069
$ws.Range("H36").Value = 0.01858977999586174
$ws.Range("I36").Value = 0.9625706352718906
$ws.Range("G37").Value = 0.02164130869719386
$ws.Range("H37").Value = 0.01768868926770477
$ws.Range("I37").Value = 0.9670921519138527
$ws.Range("G41").Value = 0.01389026083715355
$ws.Range("H41").Value = 0.01143272355649891
$ws.Range("I41").Value = 0.9887507849910595
$ws.Range("G43").Value = 0.01084397189992443
$ws.Range("H43").Value = 0.009513209651886591
$ws.Range("I43").Value = 0.9878575925218348
$ws.Range("G44").Value = 0.01257615179655825
$ws.Range("H44").Value = 0.01002227647088219
$ws.Range("I44").Value = 0.9871874643350533
